$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''42.027.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -2.47%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.218.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -4.53%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.32%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''243.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -4.02%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.620'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -3.47%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''74.29'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -2.74%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  +0.04%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.612'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -6.68%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''40.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +0.82%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0934'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -5.56%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''7.02'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -7.44%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''0.102'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -4.26%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''2.554.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -4.51%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''14.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -6.97%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''0.843'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -4.81%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''2.225.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -4.09%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''41.903.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -2.78%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''0.0₃0965'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -4.45%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''70.96'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -2.96%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''6.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -5.33%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''2.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -0.26%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''228.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -4.28%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -0.07%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''3.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -5.92%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''10.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -5.62%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''2.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -7.25%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''7.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +13.67%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  -1.91%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''167.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.15%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''20.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -4.85%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''0.0809'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -5.14%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''30.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +0.36%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.118'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -8.95%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''0.124'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -3.38%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''4.35'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -5.47%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''4.80'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.17%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.0295'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -6.66%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''13.30'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -4.57%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''2.13'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -9.54%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''5.70'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -4.02%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''110.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +4.15%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.199'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -9.82%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''59.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -5.27%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''8.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -7.07%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  -4.07%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''0.995'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.56%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''1.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -6.63%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  -3.59%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''4.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -15.23%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = '''HuobiToken'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = '''2.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -1.20%  '
$ws.Range('E51').Style = 'Normal'
